# Insert two new weekly price records for "Ají" (Terminal La Palmera de La
# Serena) right before the existing row 242, shifting the remaining rows
# down by two (so the former last two rows of the table end up duplicated
# at the new bottom of the range, rows 314-315).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 242 (this pushes old row 242 -> 244, ...,
# old row 313 -> 315, growing the sheet dimension from R313 to R315).
$ws.Rows.Item(242).Insert()
$ws.Rows.Item(242).Insert()

# --- New row 242 : Inferno / Primera -------------------------------------
$ws.Cells.Item(242, 1).Value  = 8
$ws.Cells.Item(242, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(242, 3).Value  = "Coquimbo"
$ws.Cells.Item(242, 4).Value  = 44876
$ws.Cells.Item(242, 5).Value  = 4
$ws.Cells.Item(242, 6).Value  = 100112021
$ws.Cells.Item(242, 7).Value  = "Ají"
$ws.Cells.Item(242, 8).Value  = "Inferno"
$ws.Cells.Item(242, 9).Value  = "Primera"
$ws.Cells.Item(242, 10).Value = 440
$ws.Cells.Item(242, 11).Value = 14000
$ws.Cells.Item(242, 12).Value = 15000
$ws.Cells.Item(242, 13).Value = 14500
$ws.Cells.Item(242, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(242, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(242, 16).Value = 1450
$ws.Cells.Item(242, 17).Value = 10
$ws.Cells.Item(242, 18).Value = "Hortaliza"

# --- New row 243 : Inferno / Segunda -------------------------------------
$ws.Cells.Item(243, 1).Value  = 8
$ws.Cells.Item(243, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(243, 3).Value  = "Coquimbo"
$ws.Cells.Item(243, 4).Value  = 44876
$ws.Cells.Item(243, 5).Value  = 4
$ws.Cells.Item(243, 6).Value  = 100112021
$ws.Cells.Item(243, 7).Value  = "Ají"
$ws.Cells.Item(243, 8).Value  = "Inferno"
$ws.Cells.Item(243, 9).Value  = "Segunda"
$ws.Cells.Item(243, 10).Value = 360
$ws.Cells.Item(243, 11).Value = 8000
$ws.Cells.Item(243, 12).Value = 9000
$ws.Cells.Item(243, 13).Value = 8500
$ws.Cells.Item(243, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(243, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(243, 16).Value = 850
$ws.Cells.Item(243, 17).Value = 10
$ws.Cells.Item(243, 18).Value = "Hortaliza"
